$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "E7"  = 15.31649999999999
    "C10" = -13.1389
    "C12" = -10.8906
    "E15" = 16.0997
    "C18" = -12.12949999999999
    "E20" = 15.8534
    "E29" = 16.97750000000001
    "E30" = 15.59519999999999
    "E31" = 15.9646
    "C37" = -14.11409999999999
    "E40" = 17.21710000000001
    "C55" = -13.38859999999999
    "C68" = -11.2613
    "E68" = 17.14580000000001
    "E76" = 16.12359999999999
    "C77" = -11.6604
    "C78" = -11.98850000000001
    "E87" = 16.21579999999999
    "E88" = 16.3717
    "E96" = 16.27839999999998
    "E98" = 15.43840000000001
    "E101" = 16.74650000000001
    "E102" = 16.701
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
